$d = $word.ActiveDocument

$replacements = @(
    @{old="78÷8="; new="87÷7="},
    @{old="75÷9="; new="88÷5="},
    @{old="43÷2="; new="60÷6="},
    @{old="22÷5="; new="79÷3="},
    @{old="45÷5="; new="67÷4="},
    @{old="33÷5="; new="93÷8="},
    @{old="86÷3="; new="34÷7="},
    @{old="68÷2="; new="38÷7="},
    @{old="98÷2="; new="15÷2="},
    @{old="57÷7="; new="35÷6="},
    @{old="35÷8="; new="91÷6="},
    @{old="94÷2="; new="30÷8="},
    @{old="98÷7="; new="99÷2="},
    @{old="89÷7="; new="80÷7="},
    @{old="24÷5="; new="14÷8="},
    @{old="52÷5="; new="76÷3="},
    @{old="34÷5="; new="28÷5="},
    @{old="48÷8="; new="23÷5="},
    @{old="44÷5="; new="57÷9="},
    @{old="56÷4="; new="82÷8="},
    @{old="30÷4="; new="19÷8="},
    @{old="74÷5="; new="94÷4="},
    @{old="70÷6="; new="65÷9="},
    @{old="60÷5="; new="74÷2="},
    @{old="20÷2="; new="54÷9="}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
